$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E: convert plain "duration" numbers to "Ndays" text, matching
#     the formatting already used by the rest of column E in this sheet ---
$ws.Range("E25").Value = "201days"
$ws.Range("E26").Value = "164days"
$ws.Range("E27").Value = "185days"
$ws.Range("E28").Value = "206days"
$ws.Range("E29").Value = "200days"
$ws.Range("E30").Value = "141days"
$ws.Range("E31").Value = "200days"
$ws.Range("E58").Value = "515days"
$ws.Range("E59").Value = "770days"
$ws.Range("E60").Value = "651days"
$ws.Range("E61").Value = "286days"
$ws.Range("E63").Value = "307days"
$ws.Range("E64").Value = "260days"
$ws.Range("E65").Value = "142days"
$ws.Range("E66").Value = "184days"

# --- Rows 24 & 35 (summary rows): font color changed from blue to automatic ---
$ws.Range("A24:F24").Font.ThemeColor = 1
$ws.Range("A35:F35").Font.ThemeColor = 1

# --- Update the active selection to match where the user ended up editing ---
[void]$ws.Range("E69").Select()
